$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add a new title-page paragraph "2015" right after the existing
#    "Trading Strategy Reports" title paragraph, matching its formatting
#    (centered, bold, 72-pt, cstheme minorHAnsi).
# ---------------------------------------------------------------------------

# Locate the "Trading Strategy Reports" paragraph.
$titleRange = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Trading Strategy Reports") {
        $titleRange = $p.Range
        break
    }
}

# Remember where the title paragraph sits so we can re-fetch the freshly
# inserted paragraph afterwards (indices shift once we insert).
$titleIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $titleRange.Start) {
        $titleIndex = $i
        break
    }
}

# Insert a new empty paragraph right after the title paragraph.
$endOfTitle = $titleRange.Duplicate
$endOfTitle.Collapse(0)   # wdCollapseEnd
$endOfTitle.InsertParagraphAfter()

# Copy the title paragraph's run formatting (incl. theme fonts, bold, size)
# into the freshly created paragraph, then swap its text for "2015".
$newParaRange = $d.Paragraphs($titleIndex + 1).Range
$newParaRange.FormattedText = $titleRange.FormattedText

$newParaRange = $d.Paragraphs($titleIndex + 1).Range
[void]$newParaRange.MoveEnd(1, -1)     # wdCharacter -- drop the trailing paragraph mark
$newParaRange.Text = "2015"

# ---------------------------------------------------------------------------
# 2) Update the cached result of the PAGE field in the footer from "2" to "5".
# ---------------------------------------------------------------------------

$section = $d.Sections.Item(1)
$footer = $section.Footers.Item(1)
$footerRange = $footer.Range

foreach ($fld in $footerRange.Fields) {
    if ($fld.Type -eq 33) {   # wdFieldPage
        $cachedText = $fld.Result.Text
        $count = $footerRange.Characters.Count
        for ($i = 1; $i -le $count; $i++) {
            $ch = $footerRange.Characters.Item($i)
            if ($ch.Text -eq $cachedText) {
                $ch.Text = "5"
            }
        }
    }
}
